$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 13 with the added "services" indicator, matching the formatting of the row above it
$ws.Range("D13").Value = "Объекты быт. обслу. - servicesnum (шт.) (id8001001 & 8401011)"
$ws.Range("D12").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null

# Widen column D to fit the new, longer text
$ws.Columns.Item(4).ColumnWidth = 57.33

# Move the active selection as in the authored edit
$ws.Range("D23").Select() | Out-Null
